# "Gender abbreviation"
# Replace the spelled-out gender values in the "Gender" column (column D)
# with their single-letter abbreviations:
#   male   -> M
#   female -> F
#   Female -> F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "Gender" column by inspecting the header row, falling back to
# column D (its known position in this workbook) if it cannot be found.
$genderCol = 4
$headerRange = $ws.UsedRange.Rows.Item(1)
$headerCols = $headerRange.Columns.Count
for ($c = 1; $c -le $headerCols; $c++) {
    $headerValue = $ws.Cells.Item(1, $c).Value2
    if ($headerValue -eq "Gender") {
        $genderCol = $c
    }
}

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $genderCol)
    $value = $cell.Value2
    if ($value -eq "male") {
        $cell.Value2 = "M"
        $changed = $changed + 1
    } elseif ($value -eq "female" -or $value -eq "Female") {
        $cell.Value2 = "F"
        $changed = $changed + 1
    }
}

Write-Host ("Abbreviated gender values in $changed cell(s).")
